# Applies the LOM3088.docx edit described by the commit's OOXML diff.
#
# The edit is a pure text-content rotation: the paragraph marks, run/line-break
# structure, paragraph styles and run formatting (bold labels, etc.) are all
# unchanged. Only the literal text carried by certain runs moves to a
# different paragraph. We therefore scope every Find/Replace to the owning
# paragraph's Range so the correct run is targeted unambiguously, and restrict
# MatchWildcards=$false / exact text matches throughout.

function Replace-InParagraph($ParagraphIndex, $OldText, $NewText) {
    $rng = $d.Paragraphs($ParagraphIndex).Range
    $ok = $rng.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, $NewText, 2)
    if (-not $ok) {
        throw "Find.Execute failed for paragraph $ParagraphIndex (old text not found)"
    }
}

$d = $word.ActiveDocument

# Paragraph 6 ("Objetivos" body) - single run.
Replace-InParagraph 6 "Apresentar as técnicas de caracterização e processamento de polímeros de forma que aluno seja capaz de compreender a importância de cada uma e relacioná-las com o papel desempenhado pelos polímeros na Engenharia de Materiais." "Caracterização de Polímeros: Identificação de polímeros; Determinação da massa molar média de polímeros; Determinação das propriedades físicas; análise térmica de polímeros. Processamento de materiais poliméricos"

# Paragraph 8 ("Docente(s) Responsável(eis)" body) - two runs joined by a
# manual line break (<w:br/>).
Replace-InParagraph 8 "5840897 - Clodoaldo Saron" "Apresentar as técnicas de caracterização e processamento de polímeros de forma que aluno seja capaz de compreender a importância de cada uma e relacioná-las com o papel desempenhado pelos polímeros na Engenharia de Materiais."
Replace-InParagraph 8 "1033242 - Fábio Herbst Florenzano" "Testes simples: queima e densidade relativa; FTIR de polímeros; Princípios de Ressonância Magnética Nuclear aplicada a polímeros; Viscosimetria e Reologia; Cromatografia por Exclusão de Tamanho (SEC/GPC); Termogravimetria e Calorimetria Diferencial Exploratória (DSC) de polímeros; Análise Dinâmico – Mecânica (DMA). Índice de Fluidez. Moagem e moldagem de polímeros."

# Paragraph 10 ("Programa resumido" body) - single run.
Replace-InParagraph 10 "Caracterização de Polímeros: Identificação de polímeros; Determinação da massa molar média de polímeros; Determinação das propriedades físicas; análise térmica de polímeros. Processamento de materiais poliméricos" "Experimentos desenvolvidos em laboratório didático; realização de relatórios para cada experimento."

# Paragraph 12 ("Programa" body) - single run.
Replace-InParagraph 12 "Testes simples: queima e densidade relativa; FTIR de polímeros; Princípios de Ressonância Magnética Nuclear aplicada a polímeros; Viscosimetria e Reologia; Cromatografia por Exclusão de Tamanho (SEC/GPC); Termogravimetria e Calorimetria Diferencial Exploratória (DSC) de polímeros; Análise Dinâmico – Mecânica (DMA). Índice de Fluidez. Moagem e moldagem de polímeros." "Média aritmética das notas obtidas nos relatórios. Será aprovado o aluno que obtiver nota final maior ou igual a 5,0."

# Paragraph 14 ("Avaliação" body) - six runs (bold label / value, repeated
# three times). Replace starting from the LAST run and working backwards so
# that a freshly written value can never be re-matched by a later, not yet
# processed, Find() in the same paragraph.
Replace-InParagraph 14 "Devido às características práticas da disciplina, não será oferecida recuperação." "5840897 - Clodoaldo Saron"
Replace-InParagraph 14 "Média aritmética das notas obtidas nos relatórios. Será aprovado o aluno que obtiver nota final maior ou igual a 5,0." "1. HARPER, C. A. Handbook of Plastics, Elastomers and Composites. New York: McGraw-Hill Inc, 1992 ^l2. S. V. CANEVAROLO Jr. Técnicas de Caracterização de Polímeros. São Paulo: Editora Artliber, 2005. ^l3. MANRICH, S. Processamento de Termoplásticos. Editora Artliber, 2005.^l4. NAVARRO, R.F. Fundamentos de Reologia de Polímeros. Editora da Universidade de Caxias do Sul, 1997. ^l5. MANO, E. B.; MENDES, L. C. Identificação de Plásticos, Borrachas e Fibras. Ed. Edgard Blücher, 2000. ^l6. TURI, E. A. Thermal Characterization of Polymeric Materials. New York: Academic Press, 1981.^l7. NAVARRO, R.F. Fundamentos de Reologia de Polímeros. Editora da Universidade de Caxias do Sul, 1997.MANO, E. B.; 8. MENDES, L. C. Identificação de Plásticos, Borrachas e Fibras. Ed. Edgard Blücher, 2000."
Replace-InParagraph 14 "Experimentos desenvolvidos em laboratório didático; realização de relatórios para cada experimento." "Devido às características práticas da disciplina, não será oferecida recuperação."

# Paragraph 16 ("Bibliografia" body) - a single run holding the numbered
# reference list as 7 <w:t> segments joined by manual line breaks; it
# collapses down to one plain line of text, so just overwrite Range.Text.
$d.Paragraphs(16).Range.Text = "1033242 - Fábio Herbst Florenzano"
